# Update Name of Algo
# Adjust imputed KNN values in column B to reflect the refreshed algorithm output.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value  = 7.867999999999999
$ws.Range("B7").Value  = 5.24
$ws.Range("B16").Value = 4.893000000000001
$ws.Range("B28").Value = 5.196000000000001
$ws.Range("B29").Value = 5.305999999999999
$ws.Range("B32").Value = 6.77
$ws.Range("B40").Value = 9.223000000000001
$ws.Range("B52").Value = 5.646
$ws.Range("B57").Value = 5.034000000000001
$ws.Range("B66").Value = 5.013
$ws.Range("B100").Value = 5.608
